$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# Rows containing the "Cloud - C (KRTA9AA3)" course block that must now
# reflect the merged session code KRTA9AA3/KUPT9BB1.
$courseRows = 3,4,7,8,11,12,15,16

foreach ($r in $courseRows) {
    $ws.Cells.Item($r, 1).Value = "Cloud - C (KRTA9AA3/KUPT9BB1)"
    $ws.Cells.Item($r, 2).Value = "KRTA9AA3/KUPT9BB1"
}

# Fill in the newly-specified room (column F) for each session.
# Amphi room for the 13:30 / 15:45 (week 44/47/49) and 7:45 sessions,
# room U3-4 for the second weekly block (week 45).
$amphiRows = 3,4,11,12
$u34Rows = 7,8,15,16

foreach ($r in $amphiRows) {
    $ws.Cells.Item($r, 6).Value = "U3-Amphi"
}

foreach ($r in $u34Rows) {
    $ws.Cells.Item($r, 6).Value = "U3-4"
}
